$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.336.51"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.495.20"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").Value = "2.498.20"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").Value = "67.329.16"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "2.524.36"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "2.632.14"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").Value = "  -5.90%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "516.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("E43").Value = "  -6.74%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.538"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
